$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-369) holds the "Förändrad" date. Update all entries
# from serial 45175 (2023-09-06) to 45177 (2023-09-08).
$ws.Range("C2:C369").Value = 45177
